$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"); copy H1's format (bold, border,
# centered) onto them so they share the same style as the rest of the header row.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-15 for the new columns I (I0) and J (IF)
$data = @{
    2  = @(3, 6)
    3  = @(7, 9)
    4  = @(7, 8)
    5  = @(1, 4)
    6  = @(4, 6)
    7  = @(6, 7)
    8  = @(1, 4)
    9  = @(3, 6)
    10 = @(1, 5)
    11 = @(1, 6)
    12 = @(1, 5)
    13 = @(1, 4)
    14 = @(8, 9)
    15 = @(2, 3)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
